$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.613.80"
$ws.Range("E2").Value = "  +5.57%  "
$ws.Range("D3").Value = "2.729.26"
$ws.Range("E3").Value = "  +4.48%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'590.54"
$ws.Range("E5").Value = "  +1.33%  "
$ws.Range("D6").Value = "'152.90"
$ws.Range("E6").Value = "  +6.78%  "
$ws.Range("D7").Value = "'0.995"
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("D8").Value = "'0.610"
$ws.Range("E8").Value = "  +2.14%  "
$ws.Range("D9").Value = "2.758.80"
$ws.Range("E9").Value = "  +5.25%  "
$ws.Range("D10").Value = "'6.70"
$ws.Range("E10").Value = "  +3.06%  "
$ws.Range("E11").Value = "  +7.10%  "
$ws.Range("D12").Value = "'0.391"
$ws.Range("E12").Value = "  +3.88%  "
$ws.Range("D13").Value = "'0.158"
$ws.Range("E13").Value = "  +1.74%  "
$ws.Range("D14").Value = "3.214.26"
$ws.Range("E14").Value = "  +4.60%  "
$ws.Range("D15").Value = "'26.47"
$ws.Range("E15").Value = "  +4.94%  "
$ws.Range("D16").Value = "63.465.36"
$ws.Range("E16").Value = "  +5.33%  "
$ws.Range("E17").Value = "  +8.59%  "
$ws.Range("D18").Value = "2.752.87"
$ws.Range("E18").Value = "  +5.17%  "
$ws.Range("D19").Value = "'12.03"
$ws.Range("E19").Value = "  +4.96%  "
$ws.Range("D20").Value = "'4.88"
$ws.Range("E20").Value = "  +4.22%  "
$ws.Range("D21").Value = "'365.02"
$ws.Range("E21").Value = "  +5.26%  "
$ws.Range("D22").Value = "'7.00"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("D23").Value = "'0.537"
$ws.Range("E23").Value = "  +1.12%  "
$ws.Range("D24").Value = "'0.993"
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("D25").Value = "'65.88"
$ws.Range("E25").Value = "  +3.46%  "
$ws.Range("D26").Value = "'0.167"
$ws.Range("E26").Value = "  +4.73%  "
$ws.Range("D27").Value = "'8.66"
$ws.Range("E27").Value = "  +7.93%  "
$ws.Range("D28").Value = "'0.997"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").Value = "0.0₃0905"
$ws.Range("E29").Value = "  +13.38%  "
$ws.Range("E30").Value = "  +4.45%  "
$ws.Range("D31").Value = "'7.09"
$ws.Range("E31").Value = "  +9.07%  "
$ws.Range("D32").Value = "'172.83"
$ws.Range("E32").Value = "  +2.40%  "
$ws.Range("D33").Value = "'1.20"
$ws.Range("E33").Value = "  +18.55%  "
$ws.Range("D34").Value = "'0.997"
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("D35").Value = "'20.58"
$ws.Range("E35").Value = "  +5.52%  "
$ws.Range("D36").Value = "'4.82"
$ws.Range("E36").Value = "  +12.04%  "
$ws.Range("D37").Value = "'1.43"
$ws.Range("E37").Value = "  +10.37%  "
$ws.Range("E38").Value = "  +8.99%  "
$ws.Range("E39").Value = "  +18.80%  "
$ws.Range("D40").Value = "'346.11"
$ws.Range("E40").Value = "  +8.23%  "
$ws.Range("D41").Value = "'4.22"
$ws.Range("E41").Value = "  +7.02%  "
$ws.Range("D42").Value = "'38.88"
$ws.Range("E42").Value = "  +1.00%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'5.61"
$ws.Range("E43").Value = "  +11.15%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'22.04"
$ws.Range("E44").Value = "  +10.15%  "
$ws.Range("D45").Value = "'143.15"
$ws.Range("E45").Value = "  +5.47%  "
$ws.Range("D46").Value = "'22.13"
$ws.Range("E46").Value = "  +10.38%  "
$ws.Range("E47").Value = "  +7.40%  "
$ws.Range("D48").Value = "'0.648"
$ws.Range("E48").Value = "  +6.26%  "
$ws.Range("E49").Value = "  +6.92%  "
$ws.Range("E50").Value = "  +2.81%  "
$ws.Range("D51").Value = "2.170.48"
$ws.Range("E51").Value = "  +7.20%  "
